# "minor dd lower buff" — add a new M4 lower receiver (Daniel Defense M4
# Carbine Lower Receiver) to the parts sheet, and recompute the "combo"
# rows that sum a base lower receiver with a trigger guard.
#
# Net effect on the sheet: a new row 4 is inserted for the new part, the
# combo rows are recomputed (and a second combo row added for the new part
# + trigger guard), the magwell/trigger-guard block is reshuffled, and the
# NcStar magwell row gets new recoil numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Grab the "magwell block" style (originally rows 10-13) BEFORE the
#     clean-up below touches it, so it can be re-applied once the rows it
#     belongs to land on rows 12-15 after the reshuffle. ------------------
$ws.Range("A10:P10").Copy()
$ws.Range("A21:P21").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Start from a clean slate for the data rows so nothing from the old
#     layout (stray formulas/values/styles) lingers after the reshuffle. --
$ws.Range("A3:P19").ClearFormats()
$ws.Range("A3:P19").ClearContents()

# --- Re-apply the "magwell block" style onto the rows it occupies after
#     the shuffle, rows 12-15, then discard the scratch copy. ------------
$ws.Range("A21:P21").Copy()
$ws.Range("A12:P15").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("A21:P21").Delete()

# --- Row 3: Colt M4A1 base lower — unchanged. ----------------------------
$ws.Cells.Item(3,1).Value = "colt_m4a1_std_lower_receiver"
$ws.Cells.Item(3,2).Value = "Colt M4A1"
$ws.Cells.Item(3,3).Value = 0
$ws.Cells.Item(3,4).Value = 0.3
$ws.Cells.Item(3,13).Value = 0
$ws.Cells.Item(3,14).Formula = "=C3-D3*20-E3*0.8-F3*0.6"

# --- Row 4 (new): Daniel Defense M4 Carbine Lower Receiver. -------------
$ws.Cells.Item(4,1).Value = "daniel_defense_m4_carbine_lower_receiver"
$ws.Cells.Item(4,2).Value = "Daniel Defense M4 Carbine Lower Receiver"
$ws.Cells.Item(4,3).Value = 4
$ws.Cells.Item(4,4).Value = 0.33
$ws.Cells.Item(4,13).Value = 1000
$ws.Cells.Item(4,14).Formula = "=C4-D4*20-E4*0.8-F4*0.6"

# --- Row 5: Colt M4A1 + Colt trigger guard (now row 18) combo. ----------
$ws.Cells.Item(5,2).Value = "Colt M4A1"
$ws.Cells.Item(5,3).Formula = "=C3+C18"
$ws.Cells.Item(5,4).Formula = "=D3+D18"
$ws.Cells.Item(5,5).Formula = "=E3+E18"
$ws.Cells.Item(5,6).Formula = "=F3+F18"
$ws.Cells.Item(5,14).Formula = "=C5-D5*20-E5*0.8-F5*0.6"

# --- Row 6: Daniel Defense lower + Colt trigger guard combo. ------------
$ws.Cells.Item(6,2).Value = "Daniel Defense M4 Carbine Lower Receiver"
$ws.Cells.Item(6,3).Formula = "=C4+C18"
$ws.Cells.Item(6,4).Formula = "=D4+D18"
$ws.Cells.Item(6,5).Formula = "=E4+E18"
$ws.Cells.Item(6,6).Formula = "=F4+F18"
$ws.Cells.Item(6,14).Formula = "=C6-D6*20-E6*0.8-F6*0.6"

# --- Row 7: Radian Weapons AX556 ADAC15 Lower Receiver. ------------------
$ws.Cells.Item(7,1).Value = "radian_weapons_ax556_adac15_lower_receiver"
$ws.Cells.Item(7,2).Value = "Radian Weapons AX556 ADAC15 Lower Receiver"
$ws.Cells.Item(7,3).Value = 6
$ws.Cells.Item(7,4).Value = 0.42
$ws.Cells.Item(7,5).Value = -2
$ws.Cells.Item(7,6).Value = -3
$ws.Cells.Item(7,13).Value = 4000
$ws.Cells.Item(7,14).Formula = "=C7-D7*20-E7*0.8-F7*0.6"

# --- Row 8: Noveske N4 Gen3 AR15 Lower Receiver. -------------------------
$ws.Cells.Item(8,1).Value = "noveske_n4_gen3_ar15_lower_receiver"
$ws.Cells.Item(8,2).Value = "Noveske N4 Gen3 AR15 Lower Receiver"
$ws.Cells.Item(8,3).Value = 5
$ws.Cells.Item(8,4).Value = 0.35
$ws.Cells.Item(8,5).Value = -2
$ws.Cells.Item(8,6).Value = -1
$ws.Cells.Item(8,13).Value = 2000
$ws.Cells.Item(8,14).Formula = "=C8-D8*20-E8*0.8-F8*0.6"

# --- Row 9: Aeroknox AX // 15 Lower Receiver. ----------------------------
$ws.Cells.Item(9,1).Value = "aeroknox_ax15_lower_receiver"
$ws.Cells.Item(9,2).Value = "Aeroknox AX // 15 Lower Receiver"
$ws.Cells.Item(9,3).Value = 3
$ws.Cells.Item(9,4).Value = 0.25
$ws.Cells.Item(9,5).Value = -1
$ws.Cells.Item(9,6).Value = 0
$ws.Cells.Item(9,13).Value = 1500
$ws.Cells.Item(9,14).Formula = "=C9-D9*20-E9*0.8-F9*0.6"

# --- Row 10: Fightlite SCR Rifle Lower Receiver. -------------------------
$ws.Cells.Item(10,1).Value = "fightlite_scr_rifle_lower_receiver"
$ws.Cells.Item(10,2).Value = "Fightlite SCR Rifle Lower Receiver"
$ws.Cells.Item(10,3).Value = 8
$ws.Cells.Item(10,4).Value = 0.23
$ws.Cells.Item(10,5).Value = -6
$ws.Cells.Item(10,6).Value = 12
$ws.Cells.Item(10,13).Value = 3000
$ws.Cells.Item(10,14).Formula = "=C10-D10*20-E10*0.8-F10*0.6"

# --- Row 11: spacer row (only the N formula). ----------------------------
$ws.Cells.Item(11,14).Formula = "=C11-D11*20-E11*0.8-F11*0.6"

# --- Row 12: Armaspec Rhino R-23 Tactical (magwell). ---------------------
$ws.Cells.Item(12,1).Value = "armaspec_rhino_r23_tactical_magwell"
$ws.Cells.Item(12,2).Value = "Armaspec Rhino R-23 Tactical"
$ws.Cells.Item(12,3).Value = 3
$ws.Cells.Item(12,4).Value = 0.13
$ws.Cells.Item(12,13).Value = 750
$ws.Cells.Item(12,14).Formula = "=C12-D12*20-E12*0.8-F12*0.6"

# --- Row 13: NcStar BlastAR VKARMW (magwell) — reworked numbers. ---------
$ws.Cells.Item(13,1).Value = "ncstar_blastar_vkarmw_magwell"
$ws.Cells.Item(13,2).Value = "NcStar BlastAR VKARMW"
$ws.Cells.Item(13,3).Value = 1
$ws.Cells.Item(13,4).Value = 0.07
$ws.Cells.Item(13,5).Value = -1
$ws.Cells.Item(13,6).Value = -1
$ws.Cells.Item(13,13).Value = 500
$ws.Cells.Item(13,14).Formula = "=C13-D13*20-E13*0.8-F13*0.6"

# --- Row 14: HRF Concepts RCM AR15 Milspec Magwell. ----------------------
$ws.Cells.Item(14,1).Value = "hrf_concepts_rcm_ar15_milspec_magwell"
$ws.Cells.Item(14,2).Value = "HRF Concepts RCM AR15 Milspec Magwell"
$ws.Cells.Item(14,3).Value = 2
$ws.Cells.Item(14,4).Value = 0.05
$ws.Cells.Item(14,13).Value = 250
$ws.Cells.Item(14,14).Formula = "=C14-D14*20-E14*0.8-F14*0.6"

# --- Row 15: spacer row (only the N formula, rest blank). ----------------
$ws.Cells.Item(15,14).Formula = "=C15-D15*20-E15*0.8-F15*0.6"

# --- Row 16: Magpul OD Trigger Guard. ------------------------------------
$ws.Cells.Item(16,1).Value = "magpul_trigger_guard"
$ws.Cells.Item(16,2).Value = "Magpul OD Trigger Guard"
$ws.Cells.Item(16,3).Value = 2
$ws.Cells.Item(16,4).Value = 0.02
$ws.Cells.Item(16,13).Value = 200
$ws.Cells.Item(16,14).Formula = "=C16-D16*20-E16*0.8-F16*0.6"

# --- Row 17: TI Enhanced Trigger Guard. ----------------------------------
$ws.Cells.Item(17,1).Value = "ti_enhanced_trigger_guard"
$ws.Cells.Item(17,2).Value = "TI Enhanced Trigger Guard"
$ws.Cells.Item(17,3).Value = 2.5
$ws.Cells.Item(17,4).Value = 0.05
$ws.Cells.Item(17,13).Value = 200
$ws.Cells.Item(17,14).Formula = "=C17-D17*20-E17*0.8-F17*0.6"

# --- Row 18: Colt M4 standard Trigger Guard. -----------------------------
$ws.Cells.Item(18,1).Value = "colt_m4a1_trigger_guard"
$ws.Cells.Item(18,2).Value = "Colt M4 standard Trigger Guard"
$ws.Cells.Item(18,3).Value = 2
$ws.Cells.Item(18,4).Value = 0.02
$ws.Cells.Item(18,13).Value = 0
$ws.Cells.Item(18,14).Formula = "=C18-D18*20-E18*0.8-F18*0.6"

# --- Row 19: Strike Industries Polyflex AR15 Trigger Guard. --------------
$ws.Cells.Item(19,1).Value = "strike_industries_polyflex_ar15_trigger_guard"
$ws.Cells.Item(19,2).Value = "Strike Industries Polyflex AR15 Trigger Guard"
$ws.Cells.Item(19,3).Value = 3
$ws.Cells.Item(19,4).Value = 0.08
$ws.Cells.Item(19,13).Value = 400
$ws.Cells.Item(19,14).Formula = "=C19-D19*20-E19*0.8-F19*0.6"

# --- Selection + recalc, matching the saved workbook state. -------------
$ws.Range("E4").Select()
$excel.Calculate()
